# Applies the diff: refreshes the ticker lists in columns B-F for rows 2-24,
# and appends 10 new data rows (25-34) continuing the index/ticker sequence
# in columns A-B, matching the new dimension A1:F34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New contents for existing rows 2-24 (columns B,C,D,E,F) ---
# $null marks a cell that should end up blank/empty.
$rowData = @{
    2  = @("NSE:ABSLNN50ET", "NSE:ARROWGREEN", "NSE:PHOENIXLTD", $null,           "NSE:GODREJCP")
    3  = @("NSE:ALPSINDUS",  "NSE:BLISSGVS",   $null,            $null,           "NSE:IDFCFIRSTB")
    4  = @("NSE:DIGJAMLMTD", "NSE:GHCL",       $null,            $null,           "NSE:PHOENIXLTD")
    5  = @("NSE:DLINKINDIA", "NSE:GOKEX",      $null,            $null,           $null)
    6  = @("NSE:EMAMILTD",   "NSE:HIKAL",      $null,            $null,           $null)
    7  = @("NSE:FIEMIND",    "NSE:KAUSHALYA",  $null,            $null,           $null)
    8  = @("NSE:GOCLCORP",   $null,            $null,            $null,           $null)
    9  = @("NSE:GOLDBEES",   $null,            $null,            $null,           $null)
    10 = @("NSE:GOLDETF",    $null,            $null,            $null,           $null)
    11 = @("NSE:GOLDSHARE",  $null,            $null,            $null,           $null)
    12 = @("NSE:GOYALALUM",  $null,            $null,            $null,           $null)
    13 = @("NSE:HDFCSILVER", $null,            $null,            $null,           $null)
    14 = @("NSE:IDFCFIRSTB", $null,            $null,            $null,           $null)
    15 = @("NSE:INDOBORAX",  $null,            $null,            $null,           $null)
    16 = @("NSE:INDRAMEDCO", $null,            $null,            $null,           $null)
    17 = @("NSE:IRISDOREME", $null,            $null,            $null,           $null)
    18 = @("NSE:KIRLOSBROS", $null,            $null,            $null,           $null)
    19 = @("NSE:LOTUSEYE",   $null,            $null,            $null,           $null)
    20 = @("NSE:MCL",        $null,            $null,            $null,           $null)
    21 = @("NSE:MINDTECK",   $null,            $null,            $null,           $null)
    22 = @("NSE:MONARCH",    $null,            $null,            $null,           $null)
    23 = @("NSE:MONTECARLO", $null,            $null,            $null,           $null)
    24 = @("NSE:MTARTECH",   $null,            $null,            $null,           $null)
}

$cols = @("B", "C", "D", "E", "F")

foreach ($r in 2..24) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $r)
        if ($vals[$i] -eq $null) {
            $cell.Value = ""
        } else {
            $cell.Value = $vals[$i]
        }
    }
}

# --- Append new rows 25-34 (index in A, ticker in B; C-F stay blank) ---
$newRows = @{
    25 = @(23, "NSE:NESCO")
    26 = @(24, "NSE:NFL")
    27 = @(25, "NSE:NUVOCO")
    28 = @(26, "NSE:PHOENIXLTD")
    29 = @(27, "NSE:PIDILITIND")
    30 = @(28, "NSE:POWERMECH")
    31 = @(29, "NSE:QGOLDHALF")
    32 = @(30, "NSE:RAJRATAN")
    33 = @(31, "NSE:RAMANEWS")
    34 = @(32, "NSE:RPGLIFE")
}

foreach ($r in 25..34) {
    # Copy row 24's full formatting (bold/centered/bordered index style in A,
    # plain cells in B-F) down onto the new row before writing its values.
    $ws.Range("A24:F24").Copy($ws.Range("A$r"))

    $idx = $newRows[$r][0]
    $ticker = $newRows[$r][1]

    $ws.Range("A$r").Value = $idx
    $ws.Range("B$r").Value = $ticker
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("E$r").Value = ""
    $ws.Range("F$r").Value = ""
}

Write-Output "edit applied"
